$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 63.18739966666666
$ws.Range("H2").Value = 189.562199
$ws.Range("I2").Value = 0.09596345243430386
$ws.Range("J2").Value = 0.09988075390087989
$ws.Range("Q2").Value = 5.25375847021811
$ws.Range("R2").Value = 47.283826231963
$ws.Range("S2").Value = 0.09596345243430386
$ws.Range("T2").Value = 0.09988075390087989

# Row 3
$ws.Range("I3").Value = 0.3063997713314046
$ws.Range("J3").Value = 0.3189072441572365
$ws.Range("S3").Value = 0.3063997713314046
$ws.Range("T3").Value = 0.3189072441572365

# Row 4
$ws.Range("G4").Value = 170.2928416666667
$ws.Range("H4").Value = 510.878525
$ws.Range("I4").Value = 0.2586257560429799
$ws.Range("J4").Value = 0.2691830570543736
$ws.Range("Q4").Value = 14.15911184893611
$ws.Range("R4").Value = 127.432006640425
$ws.Range("S4").Value = 0.2586257560429799
$ws.Range("T4").Value = 0.2691830570543736

# Row 5
$ws.Range("G5").Value = 77.473122
$ws.Range("H5").Value = 154.946244
$ws.Range("I5").Value = 0.1176593481802354
$ws.Range("J5").Value = 0.08164152846121862
$ws.Range("Q5").Value = 6.441554377438
$ws.Range("R5").Value = 38.649326264628
$ws.Range("S5").Value = 0.1176593481802354
$ws.Range("T5").Value = 0.08164152846121862

# Row 6
$ws.Range("G6").Value = 145.7496183333334
$ws.Range("H6").Value = 437.248855
$ws.Range("I6").Value = 0.2213516720110761
$ws.Range("J6").Value = 0.2303874164262914
$ws.Range("Q6").Value = 12.11844918273722
$ws.Range("R6").Value = 109.066042644635
$ws.Range("S6").Value = 0.2213516720110761
$ws.Range("T6").Value = 0.2303874164262914
